# New weekly price record for Coliflor (Terminal Hortofrutícola Agro Chillán)
# was inserted as a new row 162, pushing all subsequent records down by one
# row (old row 202 becomes row 203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 162; this shifts every
# existing row 162..202 down to 163..203 (and carries formatting, e.g. the
# date style on column D, down with it).
$ws.Rows(162).Insert()

# Populate the newly inserted row 162 with the new weekly record.
$ws.Cells.Item(162, 1).Value = 7
$ws.Cells.Item(162, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(162, 3).Value = "Ñuble"
$ws.Cells.Item(162, 4).Value = 44551
$ws.Cells.Item(162, 5).Value = 16
$ws.Cells.Item(162, 6).Value = 100112008
$ws.Cells.Item(162, 7).Value = "Coliflor"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 1600
$ws.Cells.Item(162, 11).Value = 800
$ws.Cells.Item(162, 12).Value = 900
$ws.Cells.Item(162, 13).Value = 850
$ws.Cells.Item(162, 14).Value = "$/unidad"
$ws.Cells.Item(162, 15).Value = "Región del Maule"
$ws.Cells.Item(162, 16).Value = 850
$ws.Cells.Item(162, 17).Value = 1
$ws.Cells.Item(162, 18).Value = "Hortaliza"
